# feat: add 2022-Q1 data
#
# 1. Insert a new worksheet "2022-Q1" right after "2021-Q4" (copying the
#    "2021-Q4" sheet's layout/styles so formatting matches the other
#    per-quarter fund sheets), and fill it in with the 2022-Q1 holdings.
# 2. Update the "总计" (totals) summary sheet with a new top row for the
#    2022-Q1 quarter, pushing the existing rows down.

$wb = $excel.ActiveWorkbook

# Helper: write $val into $rng as TEXT (not an auto-coerced number), and
# leave the cell with no explicit style - mirrors how the source data
# already stores numeric-looking figures ("37.53", "008763", ...) as
# plain inline strings with no cell style.
function Set-TextValue($rng, $val) {
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.ClearFormats()
}

# --- Step 1: add the "2022-Q1" sheet -------------------------------------
$q4 = $wb.Worksheets.Item("2021-Q4")
$q4.Copy([System.Reflection.Missing]::Value, $q4)
$newSheet = $wb.Worksheets.Item(3)
$newSheet.Name = "2022-Q1"

# Row 2 - 天弘越南市场股票（QDII）A
Set-TextValue $newSheet.Range("B2") "008763"
Set-TextValue $newSheet.Range("C2") "天弘越南市场股票（QDII）A"
Set-TextValue $newSheet.Range("D2") "37.53"
Set-TextValue $newSheet.Range("E2") "92.10"
Set-TextValue $newSheet.Range("F2") "6.59"
Set-TextValue $newSheet.Range("G2") "2.4732"
$newSheet.Range("H2").Value = 3

# Row 3 - 天弘越南市场股票（QDII）C
Set-TextValue $newSheet.Range("B3") "008764"
Set-TextValue $newSheet.Range("C3") "天弘越南市场股票（QDII）C"
Set-TextValue $newSheet.Range("D3") "14.26"
Set-TextValue $newSheet.Range("E3") "92.10"
Set-TextValue $newSheet.Range("F3") "6.59"
Set-TextValue $newSheet.Range("G3") "0.9397"
$newSheet.Range("H3").Value = 3

# --- Step 2: update the "总计" sheet --------------------------------------
$zj = $wb.Worksheets.Item("总计")

# Push the existing two rows down by one (bottom-up so we don't clobber
# data before it's been moved).
$zj.Range("A4").Value = 2
$zj.Range("B4").Value = "2021-Q3"
$zj.Range("C4").Value = 2
$zj.Range("D4").Value = 2.64

$zj.Range("A3").Value = 1
$zj.Range("B3").Value = "2021-Q4"
$zj.Range("C3").Value = 2
$zj.Range("D3").Value = 1.08

# New first data row for the 2022-Q1 quarter.
$zj.Range("A2").Value = 0
$zj.Range("B2").Value = "2022-Q1"
$zj.Range("C2").Value = 2
$zj.Range("D2").Value = 3.41

# A4 is a brand new cell - copy A3's formatting (bold/border/center style)
# onto it so the index column stays visually consistent top to bottom.
$zj.Range("A3").Copy()
$zj.Range("A4").PasteSpecial(-4122)

# Restore the originally active sheet/tab (adding + copying sheets shifts
# Excel's active-sheet selection onto the newly inserted sheet).
$wb.Worksheets.Item("2021-Q3").Activate()
